$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The employee account-statement rows (16-27) are being re-ordered:
# the "CE" / CARMEN ELENA HERRERA GARCIA block now comes first (rows 16-18,
# in ascending period order), followed by the "CC" / MIRLEY MONTALVO PEREZ
# block (rows 19-27, also in ascending period order). The Valor Mora /
# Salario Basico figures travel together with their period.

$data = @(
    @("CE", "18140237",   "CARMEN ELENA HERRERA GARCIA", "1711", 10820, 737717),
    @("CE", "18140237",   "CARMEN ELENA HERRERA GARCIA", "1712", 29509, 737717),
    @("CE", "18140237",   "CARMEN ELENA HERRERA GARCIA", "1801", 29509, 737717),
    @("CC", "1047371756", "MIRLEY MONTALVO PEREZ",       "2007", 33125, 828116),
    @("CC", "1047371756", "MIRLEY MONTALVO PEREZ",       "2008", 33125, 828116),
    @("CC", "1047371756", "MIRLEY MONTALVO PEREZ",       "2009", 33125, 828116),
    @("CC", "1047371756", "MIRLEY MONTALVO PEREZ",       "2010", 33125, 828116),
    @("CC", "1047371756", "MIRLEY MONTALVO PEREZ",       "2011", 33125, 828116),
    @("CC", "1047371756", "MIRLEY MONTALVO PEREZ",       "2012", 33125, 828116),
    @("CC", "1047371756", "MIRLEY MONTALVO PEREZ",       "2101", 33125, 828116),
    @("CC", "1047371756", "MIRLEY MONTALVO PEREZ",       "2102", 33125, 828116),
    @("CC", "1047371756", "MIRLEY MONTALVO PEREZ",       "2103", 25396, 828116)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
}
